$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.187, 0.058, -0.231, 60, 35, -44),
    @(0.206, -0.08799999999999999, 0.408, -168, -44, -56),
    @(0.467, 0.118, 0.616, -152, 68, -55),
    @(-0.046, 0.033, 0.533, -38, -35, -127),
    @(-0.294, 0.099, 0.292, -74, 47, -63),
    @(0.094, -0.107, 0.228, 113, -10, -40),
    @(0.317, 0.183, 0.082, 177, 67, 149),
    @(0.308, 0.348, 0.095, -100, -7, -106),
    @(-0.094, -0.286, 0.441, 128, 41, 61),
    @(0.078, -0.541, 0.233, 103, 59, 91)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
